# Added mean curve method1
# Adds a new "exclude" column (E) with per-experiment exclusion lists.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("E1").Value = "exclude"

# Values for the "exclude" column, one per experiment row (2-21)
$excludeValues = @(
    "außen,3l_kü,3l_wz,4l,2l,1l,bd,bd_original",  # row 2
    "außen,3l_kü,3l_wz,4l,2l,1l,bd,bd_original",  # row 3
    "außen,3l_kü,3l_wz,4l,2l,1l,bd,bd_original",  # row 4
    "außen,3l_kü,3l_wz,4l,2l,1l,bd,bd_original",  # row 5
    "außen,3l_kü,3l_wz,4l,2l,1l,bd,bd_original",  # row 6
    "außen,3l_kü,3l_wz,4l,2l,1l,bd,bd_original",  # row 7
    "außen,3l_kü,3l_wz,4l,2l,1l",                 # row 8
    "außen,3l_kü,3l_wz,4l,2l,1l",                 # row 9
    "außen,3l_kü,3l_wz,4l,2l,1l",                 # row 10
    "außen,3l_kü,3l_wz,4l,2l,1l",                 # row 11
    "außen,3l_kü,3l_wz,4l,2l,1l",                 # row 12
    "außen,3l_kü,3l_wz,4l,2l,1l",                 # row 13
    "1a,1l,2l,3l,4a_testo,4l,außen",              # row 14
    "1a,1l,2l,3l,4a_testo,4l,außen",              # row 15
    "1a,1l,2l,3l,4a_testo,4l,außen",              # row 16
    "1l,1a_testo,1t,2a_testo,2l,3l,4a_testo,4l,außen",  # row 17
    "1a,1l,1l_sub,2l,3f,3l,4a,4a_sub,4a_testo,4l,5flur,5treppe,außen,tr,weather",  # row 18
    "1a,1l,1l_sub,2l,3f,3l,4a,4a_sub,4a_testo,4l,5flur,5treppe,außen,tr,weather",  # row 19
    "1a,1l,1l_sub,2l,3f,3l,4a,4a_sub,4a_testo,4l,5flur,5treppe,außen,tr,weather",  # row 20
    "1a,1l,1l_sub,2l,3f,3l,4a,4a_sub,4a_testo,4l,5flur,5treppe,außen,tr,weather"   # row 21
)

for ($i = 0; $i -lt $excludeValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $excludeValues[$i]
}

# Column E width, matching the source layout (stored width rounds to 68)
$ws.Columns.Item(5).ColumnWidth = 67.14286

# Update the selected cell to match the new active selection
$ws.Range("E9").Select()
